# Adds a new "2022-Q3" sheet (fund holdings) positioned between "总计" and
# "2022-Q2", and updates the "总计" summary sheet with a new row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before the existing "2022-Q2"
#    sheet (so the tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1).
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$q3rows = @(
    @("310328", "申万菱信新动力混合",                 "27.57", "90.17", "3.06", "0.8436", 9),
    @("008188", "前海开源稳健增长三年持有期混合",       "21.72", "66.43", "2.40", "0.5213", 9),
    @("010887", "南方消费升级混合A",                   "8.90",  "69.72", "4.08", "0.3631", 5),
    @("010888", "南方消费升级混合C",                   "3.45",  "69.72", "4.08", "0.1408", 5),
    @("310388", "申万菱信消费增长混合A",                "2.87",  "91.17", "3.72", "0.1068", 10),
    @("011287", "前海开源聚慧三年持有期混合",           "2.84",  "66.61", "2.38", "0.0676", 9),
    @("006775", "前海开源优质成长混合",                 "2.48",  "69.80", "2.49", "0.0618", 10),
    @("002293", "南方益和灵活配置混合",                 "1.12",  "75.22", "2.80", "0.0314", 7),
    @("002407", "前海开源恒远灵活配置混合",             "1.03",  "67.57", "2.50", "0.0258", 9),
    @("006216", "前海开源价值成长灵活配置混合A",        "1.11",  "65.31", "2.28", "0.0253", 10),
    @("006217", "前海开源价值成长灵活配置混合C",        "0.47",  "65.31", "2.28", "0.0107", 10),
    @("015254", "申万菱信消费增长混合C",                "0.05",  "91.17", "3.72", "0.0019", 10)
)

# Columns B..G on the data rows hold numeric-looking values that must stay
# TEXT (matching the other sheets, which store them as inlineStr/strings),
# so force a text number-format before writing them.
$wsQ3.Range("B1:G13").NumberFormat = "@"

# Header row (B1:H1) — bold, bordered, centered, matching the other sheets.
$headerRng = $wsQ3.Range("B1:H1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

$col = 2
foreach ($h in $headers) {
    $wsQ3.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Column A (row index) — same bold/border/centered look as the header.
$idRng = $wsQ3.Range("A2:A13")
$idRng.Font.Bold = $true
$idRng.HorizontalAlignment = -4108
$idRng.VerticalAlignment = -4160
$idRng.Borders.LineStyle = 1

$r = 2
$idx = 0
foreach ($row in $q3rows) {
    $wsQ3.Cells.Item($r, 1).Value = $idx
    $wsQ3.Cells.Item($r, 2).Value = $row[0]
    $wsQ3.Cells.Item($r, 3).Value = $row[1]
    $wsQ3.Cells.Item($r, 4).Value = $row[2]
    $wsQ3.Cells.Item($r, 5).Value = $row[3]
    $wsQ3.Cells.Item($r, 6).Value = $row[4]
    $wsQ3.Cells.Item($r, 7).Value = $row[5]
    $wsQ3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
    $idx = $idx + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row 2 for "2022-Q3" (shifting the
#    existing "2022-Q2"/"2022-Q1" rows down) and renumber the A-column index.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 12
$wsTotal.Cells.Item(2, 4).Value = 2.2

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 18
$wsTotal.Cells.Item(3, 4).Value = 3.33

$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(4, 3).Value = 25
$wsTotal.Cells.Item(4, 4).Value = 5.84

# Keep "总计" as the active/selected tab (matches the unmodified bookViews
# state — only the sheet list itself changes per the diff).
$wsTotal.Activate()
